$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename G1/H1, add new header I1 ---
$ws.Range("G1").Value = "battery_size_70000"
$ws.Range("H1").Value = "battery_size_110000"

# Clone formatting of G1 (bold, border, centered) onto the new I1 header cell
$ws.Range("G1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "battery_size_150000"
$excel.CutCopyMode = $false

# --- Data rows: updated G/H values, new I column values ---
$ws.Range("G2").Value = 1103.722925000901
$ws.Range("H2").Value = 1734.42173928713
$ws.Range("I2").Value = 2365.120553573359
$ws.Range("G3").Value = 1446.328804257992
$ws.Range("H3").Value = 2272.80240669113
$ws.Range("I3").Value = 3099.276009124269
$ws.Range("G4").Value = 2190.710524043857
$ws.Range("H4").Value = 3442.545109211776
$ws.Range("I4").Value = 4694.379694379693
$ws.Range("G5").Value = 2114.583333333333
$ws.Range("H5").Value = 3322.916666666667
$ws.Range("I5").Value = 4531.25
$ws.Range("G6").Value = 2159.942117807286
$ws.Range("H6").Value = 3394.194756554307
$ws.Range("I6").Value = 4628.447395301328
$ws.Range("G7").Value = 3426.115083418454
$ws.Range("H7").Value = 5383.895131086142
$ws.Range("I7").Value = 7341.675178753831
$ws.Range("G8").Value = 3551.261690488795
$ws.Range("H8").Value = 5580.55408505382
$ws.Range("I8").Value = 7609.846479618847
$ws.Range("G9").Value = 6979.555461435971
$ws.Range("H9").Value = 10967.87286797081
$ws.Range("I9").Value = 14956.19027450565
$ws.Range("G10").Value = 8912.249963498321
$ws.Range("H10").Value = 14004.9642283545
$ws.Range("I10").Value = 19097.67849321069
$ws.Range("G11").Value = 8948.894470803485
$ws.Range("H11").Value = 14062.54845411976
$ws.Range("I11").Value = 19176.20243743604
$ws.Range("G12").Value = 9758.4361118009
$ws.Range("H12").Value = 15334.68531854427
$ws.Range("I12").Value = 20910.93452528764
$ws.Range("G13").Value = 11167.55258344764
$ws.Range("H13").Value = 17549.01120256058
$ws.Range("I13").Value = 23930.46982167352
$ws.Range("G14").Value = 499.9857146938659
$ws.Range("H14").Value = 785.691837376075
$ws.Range("I14").Value = 1071.397960058284
$ws.Range("G15").Value = 700.0700070007
$ws.Range("H15").Value = 1100.1100110011
$ws.Range("I15").Value = 1500.1500150015
$ws.Range("G16").Value = 795.3822379783657
$ws.Range("H16").Value = 1249.886373966003
$ws.Range("I16").Value = 1704.390509953641
$ws.Range("G17").Value = 1955.767307805423
$ws.Range("H17").Value = 3073.348626551379
$ws.Range("I17").Value = 4190.929945297336
$ws.Range("G18").Value = 3914.348063284233
$ws.Range("H18").Value = 6151.118385160938
$ws.Range("I18").Value = 8387.888707037642
$ws.Range("G19").Value = 6771.370526644788
$ws.Range("H19").Value = 10640.72511329895
$ws.Range("I19").Value = 14510.07969995312
$ws.Range("G20").Value = 6616.883513530022
$ws.Range("H20").Value = 10397.95980697575
$ws.Range("I20").Value = 14179.03610042148
$ws.Range("G21").Value = 7879.614475359155
$ws.Range("H21").Value = 12382.25131842153
$ws.Range("I21").Value = 16884.8881614839
$ws.Range("G22").Value = 8497.078087618151
$ws.Range("H22").Value = 13352.55128054281
$ws.Range("I22").Value = 18208.02447346747
$ws.Range("G23").Value = 9025.955088947214
$ws.Range("H23").Value = 14183.64371120277
$ws.Range("I23").Value = 19341.33233345831
$ws.Range("G24").Value = 9025.955088947214
$ws.Range("H24").Value = 14183.64371120277
$ws.Range("I24").Value = 19341.33233345831
$ws.Range("G25").Value = 1038.065215917565
$ws.Range("H25").Value = 1631.245339299031
$ws.Range("I25").Value = 2224.425462680496
$ws.Range("G26").Value = 1134.259259259259
$ws.Range("H26").Value = 1782.407407407407
$ws.Range("I26").Value = 2430.555555555556
$ws.Range("G27").Value = 3744.534375290724
$ws.Range("H27").Value = 5884.268304028282
$ws.Range("I27").Value = 8024.002232765838
$ws.Range("G28").Value = 8525.641025641024
$ws.Range("H28").Value = 13397.43589743589
$ws.Range("I28").Value = 18269.23076923077
$ws.Range("G29").Value = 5194.218608852754
$ws.Range("H29").Value = 8162.343528197185
$ws.Range("I29").Value = 11130.46844754162
$ws.Range("G30").Value = 11876.69667095299
$ws.Range("H30").Value = 18663.38048292613
$ws.Range("I30").Value = 25450.06429489927
$ws.Range("G31").Value = 9236.111111111111
$ws.Range("H31").Value = 14513.88888888889
$ws.Range("I31").Value = 19791.66666666667
$ws.Range("G32").Value = 7794.186591654946
$ws.Range("H32").Value = 12248.00750117206
$ws.Range("I32").Value = 16701.82841068917
$ws.Range("G33").Value = 7794.186591654946
$ws.Range("H33").Value = 12248.00750117206
$ws.Range("I33").Value = 16701.82841068917
